$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 1.988074333333333
$ws.Range("H2").Value = 5.964223
$ws.Range("I2").Value = 0.01657769708907969
$ws.Range("J2").Value = 0.01657769708907968
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 334.2148329024984
$ws.Range("R2").Value = 3007.933496122486
$ws.Range("S2").Value = 0.004947086953226895
$ws.Range("T2").Value = 0.004947086953226895
$ws.Range("G3").Value = 1.988074333333333
$ws.Range("H3").Value = 5.964223
$ws.Range("I3").Value = 0.01657769708907969
$ws.Range("J3").Value = 0.01657769708907968
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("Q3").Value = 324.0685133021846
$ws.Range("R3").Value = 2916.616619719661
$ws.Range("S3").Value = 0.004796899946617806
$ws.Range("T3").Value = 0.004796899946617805
$ws.Range("G4").Value = 1.988074333333333
$ws.Range("H4").Value = 5.964223
$ws.Range("I4").Value = 0.01657769708907969
$ws.Range("J4").Value = 0.01657769708907968
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 330.0074764923966
$ws.Range("R4").Value = 2970.06728843157
$ws.Range("S4").Value = 0.004884809172724968
$ws.Range("T4").Value = 0.004884809172724967
$ws.Range("G5").Value = 1.988074333333333
$ws.Range("H5").Value = 5.964223
$ws.Range("I5").Value = 0.01657769708907969
$ws.Range("J5").Value = 0.01657769708907968
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 131.6636707085854
$ws.Range("R5").Value = 1184.973036377269
$ws.Range("S5").Value = 0.001948901016510019
$ws.Range("T5").Value = 0.001948901016510018
$ws.Range("I6").Value = 0.7746030815641455
$ws.Range("J6").Value = 0.7746030815641454
$ws.Range("M6").Value = 168.1098273333333
$ws.Range("N6").Value = 504.329482
$ws.Range("O6").Value = 0.2984182258032519
$ws.Range("P6").Value = 0.298418225803252
$ws.Range("Q6").Value = 15616.39340371692
$ws.Range("R6").Value = 140547.5406334523
$ws.Range("S6").Value = 0.231155677302104
$ws.Range("T6").Value = 0.231155677302104
$ws.Range("I7").Value = 0.7746030815641455
$ws.Range("J7").Value = 0.7746030815641454
$ws.Range("O7").Value = 0.2893586437755394
$ws.Range("P7").Value = 0.2893586437755394
$ws.Range("S7").Value = 0.2241380971457547
$ws.Range("T7").Value = 0.2241380971457547
$ws.Range("I8").Value = 0.7746030815641455
$ws.Range("J8").Value = 0.7746030815641454
$ws.Range("M8").Value = 165.99353
$ws.Range("N8").Value = 497.98059
$ws.Range("O8").Value = 0.294661504941043
$ws.Range("P8").Value = 0.294661504941043
$ws.Range("Q8").Value = 15419.80209052117
$ws.Range("R8").Value = 138778.2188146906
$ws.Range("S8").Value = 0.2282457097456606
$ws.Range("T8").Value = 0.2282457097456606
$ws.Range("I9").Value = 0.7746030815641455
$ws.Range("J9").Value = 0.7746030815641454
$ws.Range("M9").Value = 66.22673433333334
$ws.Range("N9").Value = 198.680203
$ws.Range("O9").Value = 0.1175616254801657
$ws.Range("P9").Value = 0.1175616254801657
$ws.Range("Q9").Value = 6152.065905951418
$ws.Range("R9").Value = 55368.59315356275
$ws.Range("S9").Value = 0.09106359737062629
$ws.Range("T9").Value = 0.09106359737062628
$ws.Range("G10").Value = 23.741365
$ws.Range("H10").Value = 71.224095
$ws.Range("I10").Value = 0.1979690350870239
$ws.Range("J10").Value = 0.1979690350870239
$ws.Range("M10").Value = 168.1098273333333
$ws.Range("N10").Value = 504.329482
$ws.Range("O10").Value = 0.2984182258032519
$ws.Range("P10").Value = 0.298418225803252
$ws.Range("Q10").Value = 3991.156770807644
$ws.Range("R10").Value = 35920.41093726879
$ws.Range("S10").Value = 0.05907756821465141
$ws.Range("T10").Value = 0.05907756821465143
$ws.Range("G11").Value = 23.741365
$ws.Range("H11").Value = 71.224095
$ws.Range("I11").Value = 0.1979690350870239
$ws.Range("J11").Value = 0.1979690350870239
$ws.Range("O11").Value = 0.2893586437755394
$ws.Range("P11").Value = 0.2893586437755394
$ws.Range("Q11").Value = 3869.990538238352
$ws.Range("R11").Value = 34829.91484414517
$ws.Range("S11").Value = 0.05728405150233343
$ws.Range("T11").Value = 0.05728405150233343
$ws.Range("G12").Value = 23.741365
$ws.Range("H12").Value = 71.224095
$ws.Range("I12").Value = 0.1979690350870239
$ws.Range("J12").Value = 0.1979690350870239
$ws.Range("M12").Value = 165.99353
$ws.Range("N12").Value = 497.98059
$ws.Range("O12").Value = 0.294661504941043
$ws.Range("P12").Value = 0.294661504941043
$ws.Range("Q12").Value = 3940.91298336845
$ws.Range("R12").Value = 35468.21685031605
$ws.Range("S12").Value = 0.05833385381046861
$ws.Range("T12").Value = 0.05833385381046861
$ws.Range("G13").Value = 23.741365
$ws.Range("H13").Value = 71.224095
$ws.Range("I13").Value = 0.1979690350870239
$ws.Range("J13").Value = 0.1979690350870239
$ws.Range("M13").Value = 66.22673433333334
$ws.Range("N13").Value = 198.680203
$ws.Range("O13").Value = 0.1175616254801657
$ws.Range("P13").Value = 0.1175616254801657
$ws.Range("Q13").Value = 1572.313072565699
$ws.Range("R13").Value = 14150.81765309129
$ws.Range("S13").Value = 0.02327356155957048
$ws.Range("T13").Value = 0.02327356155957048
$ws.Range("G14").Value = 1.301204666666667
$ws.Range("H14").Value = 3.903614
$ws.Range("I14").Value = 0.01085018625975097
$ws.Range("J14").Value = 0.01085018625975097
$ws.Range("M14").Value = 168.1098273333333
$ws.Range("N14").Value = 504.329482
$ws.Range("O14").Value = 0.2984182258032519
$ws.Range("P14").Value = 0.298418225803252
$ws.Range("Q14").Value = 218.7452918386609
$ws.Range("R14").Value = 1968.707626547948
$ws.Range("S14").Value = 0.003237893333269707
$ws.Range("T14").Value = 0.003237893333269707
$ws.Range("G15").Value = 1.301204666666667
$ws.Range("H15").Value = 3.903614
$ws.Range("I15").Value = 0.01085018625975097
$ws.Range("J15").Value = 0.01085018625975097
$ws.Range("O15").Value = 0.2893586437755394
$ws.Range("P15").Value = 0.2893586437755394
$ws.Range("Q15").Value = 212.1044745452331
$ws.Range("R15").Value = 1908.940270907098
$ws.Range("S15").Value = 0.003139595180833534
$ws.Range("T15").Value = 0.003139595180833534
$ws.Range("G16").Value = 1.301204666666667
$ws.Range("H16").Value = 3.903614
$ws.Range("I16").Value = 0.01085018625975097
$ws.Range("J16").Value = 0.01085018625975097
$ws.Range("M16").Value = 165.99353
$ws.Range("N16").Value = 497.98059
$ws.Range("O16").Value = 0.294661504941043
$ws.Range("P16").Value = 0.294661504941043
$ws.Range("Q16").Value = 215.9915558724734
$ws.Range("R16").Value = 1943.92400285226
$ws.Range("S16").Value = 0.003197132212188848
$ws.Range("T16").Value = 0.003197132212188847
$ws.Range("G17").Value = 1.301204666666667
$ws.Range("H17").Value = 3.903614
$ws.Range("I17").Value = 0.01085018625975097
$ws.Range("J17").Value = 0.01085018625975097
$ws.Range("M17").Value = 66.22673433333334
$ws.Range("N17").Value = 198.680203
$ws.Range("O17").Value = 0.1175616254801657
$ws.Range("P17").Value = 0.1175616254801657
$ws.Range("Q17").Value = 86.17453577262691
$ws.Range("R17").Value = 775.570821953642
$ws.Range("S17").Value = 0.001275565533458883
$ws.Range("T17").Value = 0.001275565533458883
